$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.194888178913738
$ws.Range("C2").Value = 0.5878594249201278
$ws.Range("J2").Value = 0.009584664536741214
$ws.Range("P2").Value = 0.134185303514377
$ws.Range("S2").Value = 0.07348242811501597
$ws.Range("C3").Value = 0.01075268817204301
$ws.Range("J3").Value = 0.04838709677419355
$ws.Range("P3").Value = 0.7580645161290323
$ws.Range("S3").Value = 0.1827956989247312
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.05670103092783505
$ws.Range("D6").Value = 0.01030927835051546
$ws.Range("F6").Value = 0.02061855670103093
$ws.Range("J6").Value = 0.2371134020618557
$ws.Range("O6").Value = 0.0154639175257732
$ws.Range("Q6").Value = 0.1649484536082474
$ws.Range("R6").Value = 0.08247422680412371
$ws.Range("S6").Value = 0.4123711340206185
$ws.Range("B7").Value = 0.195
$ws.Range("D7").Value = 0.035
$ws.Range("F7").Value = 0.02
$ws.Range("J7").Value = 0.13
$ws.Range("O7").Value = 0.005
$ws.Range("Q7").Value = 0.165
$ws.Range("R7").Value = 0.07000000000000001
$ws.Range("S7").Value = 0.38
$ws.Range("B8").Value = 0.112
$ws.Range("D8").Value = 0.01333333333333333
$ws.Range("F8").Value = 0.07733333333333334
$ws.Range("J8").Value = 0.1173333333333333
$ws.Range("O8").Value = 0.01066666666666667
$ws.Range("Q8").Value = 0.1706666666666667
$ws.Range("R8").Value = 0.09066666666666667
$ws.Range("S8").Value = 0.408
$ws.Range("B9").Value = 0.1274509803921569
$ws.Range("D9").Value = 0.004901960784313725
$ws.Range("F9").Value = 0.06372549019607843
$ws.Range("J9").Value = 0.1274509803921569
$ws.Range("O9").Value = 0.02941176470588235
$ws.Range("Q9").Value = 0.142156862745098
$ws.Range("R9").Value = 0.09803921568627451
$ws.Range("S9").Value = 0.4068627450980392
$ws.Range("B10").Value = 0.1056845476381105
$ws.Range("D10").Value = 0.02081665332265813
$ws.Range("E10").Value = 0.0008006405124099279
$ws.Range("F10").Value = 0.06725380304243395
$ws.Range("J10").Value = 0.1281024819855885
$ws.Range("O10").Value = 0.01040832666132906
$ws.Range("Q10").Value = 0.211369095276221
$ws.Range("R10").Value = 0.08086469175340272
$ws.Range("S10").Value = 0.3746997598078463
$ws.Range("G11").Value = 0.1461794019933555
$ws.Range("J11").Value = 0.07973421926910298
$ws.Range("K11").Value = 0.1993355481727575
$ws.Range("L11").Value = 0.5548172757475083
$ws.Range("S11").Value = 0.01993355481727575
$ws.Range("G12").Value = 0.7894736842105263
$ws.Range("J12").Value = 0.1637426900584795
$ws.Range("K12").Value = 0.005847953216374269
$ws.Range("L12").Value = 0.02923976608187134
$ws.Range("S12").Value = 0.01169590643274854
$ws.Range("G13").Value = 0.5909090909090909
$ws.Range("J13").Value = 0.3181818181818182
$ws.Range("S13").Value = 0.09090909090909091
$ws.Range("F15").Value = 0.02173913043478261
$ws.Range("H15").Value = 0.07608695652173914
$ws.Range("I15").Value = 0.08152173913043478
$ws.Range("J15").Value = 0.391304347826087
$ws.Range("K15").Value = 0.09239130434782608
$ws.Range("M15").Value = 0.03260869565217391
$ws.Range("O15").Value = 0.03804347826086957
$ws.Range("S15").Value = 0.266304347826087
$ws.Range("F16").Value = 0.0101010101010101
$ws.Range("H16").Value = 0.1363636363636364
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3939393939393939
$ws.Range("K16").Value = 0.1313131313131313
$ws.Range("M16").Value = 0.0202020202020202
$ws.Range("O16").Value = 0.06060606060606061
$ws.Range("S16").Value = 0.1565656565656566
$ws.Range("F17").Value = 0.01425178147268409
$ws.Range("H17").Value = 0.171021377672209
$ws.Range("I17").Value = 0.0997624703087886
$ws.Range("J17").Value = 0.3871733966745843
$ws.Range("K17").Value = 0.1092636579572447
$ws.Range("M17").Value = 0.009501187648456057
$ws.Range("N17").Value = 0.002375296912114014
$ws.Range("O17").Value = 0.07125890736342043
$ws.Range("S17").Value = 0.1353919239904988
$ws.Range("F18").Value = 0.0108695652173913
$ws.Range("H18").Value = 0.1521739130434783
$ws.Range("I18").Value = 0.1032608695652174
$ws.Range("J18").Value = 0.4402173913043478
$ws.Range("K18").Value = 0.1141304347826087
$ws.Range("M18").Value = 0.005434782608695652
$ws.Range("O18").Value = 0.04347826086956522
$ws.Range("S18").Value = 0.1304347826086956
$ws.Range("F19").Value = 0.01430842607313196
$ws.Range("H19").Value = 0.1868044515103339
$ws.Range("I19").Value = 0.08823529411764706
$ws.Range("J19").Value = 0.3863275039745628
$ws.Range("K19").Value = 0.1009538950715421
$ws.Range("M19").Value = 0.02305246422893482
$ws.Range("O19").Value = 0.05961844197138315
$ws.Range("S19").Value = 0.1383147853736089
